# Suppress any confirmation dialogs (sheet deletion, etc.)
$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add a new "input" row (single cell, value 0, column A) at the bottom of
#    each of the 29 surviving course worksheets ("Add input files for all
#    courses").
# ---------------------------------------------------------------------------
$newRows = @{
    1  = 17
    2  = 17
    3  = 17
    4  = 19
    5  = 19
    6  = 18
    7  = 15
    8  = 19
    9  = 20
    10 = 20
    11 = 22
    12 = 20
    13 = 16
    14 = 14
    15 = 22
    16 = 20
    17 = 20
    18 = 16
    19 = 17
    20 = 13
    21 = 23
    22 = 12
    23 = 17
    24 = 12
    25 = 16
    26 = 17
    27 = 20
    28 = 19
    29 = 11
}

foreach ($idx in 1..29) {
    $ws = $wb.Worksheets.Item($idx)
    $r = $newRows[$idx]
    $ws.Cells.Item($r, 1).Value = 0
}

# Sheet 10 ("13cen") also had its existing row-19 "Nr." cell (A19, previously
# the shared string "18") re-entered as a literal number 18.
$ws10 = $wb.Worksheets.Item(10)
$ws10.Cells.Item(19, 1).Value = 18

# Put the selection on the freshly-added row of each sheet, matching the
# per-sheet <selection> the workbook was saved with.
foreach ($idx in 1..29) {
    $ws = $wb.Worksheets.Item($idx)
    $r = $newRows[$idx]
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 1)).Select()
}

# ---------------------------------------------------------------------------
# 2. Drop the six "Sport" elective sheets - they're no longer part of the
#    course list.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("13spo23").Delete()
$wb.Worksheets.Item("13spo22").Delete()
$wb.Worksheets.Item("13spo21").Delete()
$wb.Worksheets.Item("13spo13").Delete()
$wb.Worksheets.Item("13spo12").Delete()
$wb.Worksheets.Item("13spo11").Delete()

# ---------------------------------------------------------------------------
# 3. The active tab ends up back on "13spaF" (the old last sheet before the
#    deleted block), which is also where Excel naturally lands the selection
#    after the deletions above - select it explicitly to be sure.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("13spaF").Select()
